$d = $word.ActiveDocument

# Locate the "Version: ... " / "Date: ..." paragraph robustly (rather
# than assuming a fixed paragraph index) so the edit still lands
# correctly even if the document gains/loses paragraphs elsewhere.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Version:*") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    # Each Find/Replace below re-reads $target.Range fresh (instead of
    # reusing one Range variable) because Find.Execute collapses its
    # Range to the last match; reusing it would let the next search
    # wander outside the paragraph and corrupt unrelated text (e.g. the
    # "board=53.0" support URL also contains the digit "3").

    # Version: 4.15 -> 4.16  (the lone digit run right after "4.1")
    $target.Range.Find.Execute("5", $true, $false, $false, $false, $false,
                                $true, 1, $false, "6", 2) | Out-Null

    # Date: 6/3/2020 -> 6/22/2020  (the lone digit run right after "6/")
    $target.Range.Find.Execute("3", $true, $false, $false, $false, $false,
                                $true, 1, $false, "22", 2) | Out-Null
}
